$d = $word.ActiveDocument

# The paragraph that currently contains "...впорядкування та доступн" + bookmark + "ої інформації."
# is replaced wholesale with the expanded content: the merged sentence, the new
# "Розподіл підзадач" section (with sub-contributors and their task lists), the
# "План роботи за стандартом ISO/IEC 12207" section, and a final paragraph that
# keeps the _GoBack bookmark.
$newXml = @'
<w:p w:rsidR="001C40B1" w:rsidRDefault="001C40B1" w:rsidP="001C40B1">
      <w:pPr>
        <w:pStyle w:val="a9"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="37"/>
        </w:numPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:left="0" w:firstLine="709"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Покращення</w:t>
      </w:r>
      <w:r w:rsidRPr="00C05A8D">
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve"> впорядкування та доступної інформації</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a9"/>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:left="709"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="709"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Розподіл підзадач</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="709"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Дмитренко Ярослава</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a9"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="38"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="1134"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:left="0" w:firstLine="709"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>В</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>и</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>ведення</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve"> всього реєстру на екран.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a9"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="38"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="1134"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:left="0" w:firstLine="709"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>В</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>и</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>ведення</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve"> всього реєстру в заданий текстовий файл.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a9"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="38"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="1134"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:left="0" w:firstLine="709"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Додавання нового запису до реєстру.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="709"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Середа Марина</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a9"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="39"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="1134"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:left="0" w:firstLine="709"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>П</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>ошук запису в реєстрі за заданим державним номером (якщо запис відсутній, виводиться відповідне повідомлення)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a9"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="39"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="1134"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:left="0" w:firstLine="709"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>В</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>илучення заданого запису з реєстру</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a9"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="1134"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:left="0" w:firstLine="709"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Безкровна Вероніка</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a9"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="40"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="1134"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:left="0" w:firstLine="709"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>З</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>авершення роботи програми з автоматичним записом реєстру у файл.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a9"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="40"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="1134"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:left="0" w:firstLine="709"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Автоматичне завантаження реєстру з файлу під час запуску програми.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="1134"/>
        </w:tabs>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="360" w:lineRule="auto"/>
        <w:ind w:firstLine="709"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve">План роботи </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="ru-RU"/>
        </w:rPr>
        <w:t>за стандартом ISO/IEC 12207</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a9"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="41"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="1134"/>
          <w:tab w:val="left" w:pos="1276"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:left="0" w:firstLine="709"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Аналіз задач.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a9"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="41"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="1134"/>
          <w:tab w:val="left" w:pos="1276"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:left="0" w:firstLine="709"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Розподіл функцій.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a9"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="41"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="1134"/>
          <w:tab w:val="left" w:pos="1276"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:left="0" w:firstLine="709"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Планування</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve"> роботи.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a9"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="41"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="1134"/>
          <w:tab w:val="left" w:pos="1276"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:left="0" w:firstLine="709"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Реалізація модулів</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a9"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="41"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="1134"/>
          <w:tab w:val="left" w:pos="1276"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:left="0" w:firstLine="709"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Інтеграці</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>я модулів.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a9"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="41"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="1134"/>
          <w:tab w:val="left" w:pos="1276"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:ind w:left="0" w:firstLine="709"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Верифікація відповідності вимогам</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
'@

$targetPara = $d.Paragraphs.Last
$targetRange = $targetPara.Range
Write-Host "Target paragraph text before edit:" $targetRange.Text
$targetRange.InsertXML($newXml)

Write-Host "Paragraphs after edit:" $d.Paragraphs.Count
